$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF2:BF31 hold a "Date" column that was mis-derived as "11-28-2013-14"
# (day concatenated with season string). Fix it to the correct ISO date
# string "2013-11-28" for every data row, keeping it as literal text
# (not an Excel date serial) and without altering the cell's formatting.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    $cell.NumberFormat = "@"
    $cell.Value = "2013-11-28"
    $cell.ClearFormats()
}
